$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 88
$ws.Range("H88").Value = 31144.4
$ws.Range("I88").Value = 775
$ws.Range("J88").Value = 38736.75
$ws.Range("K88").Value = 775
$ws.Range("L88").Value = 38736.75
$ws.Range("M88").Value = -369
$ws.Range("N88").Value = -39548.75
# Row 91
$ws.Range("H91").Value = 31144.4
$ws.Range("I91").Value = 775
$ws.Range("J91").Value = 38736.75
$ws.Range("K91").Value = 775
$ws.Range("L91").Value = 38736.75
$ws.Range("M91").Value = 629
$ws.Range("N91").Value = -41544.75
# Row 97
$ws.Range("H97").Value = 3001.4285
$ws.Range("J97").Value = 3001.4285
$ws.Range("L97").Value = 9004.2855
$ws.Range("N97").Value = -9996.2855
# Row 100
$ws.Range("H100").Value = 2490.8333
$ws.Range("I100").Value = 2440
$ws.Range("J100").Value = 2541.6667
$ws.Range("K100").Value = 2440
$ws.Range("L100").Value = 2541.6667
$ws.Range("M100").Value = -1899
$ws.Range("N100").Value = -3623.6667
# Row 103
$ws.Range("H103").Value = 143356.14
$ws.Range("I103").Value = 143356.14
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 430068.42
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -429482.42
$ws.Range("N103").ClearContents() | Out-Null
# Row 106
$ws.Range("H106").Value = 14289172
$ws.Range("I106").Value = 20003200
$ws.Range("J106").Value = 4103
$ws.Range("K106").Value = 20003200
$ws.Range("L106").Value = 4103
$ws.Range("M106").Value = -20002569
$ws.Range("N106").Value = -5365
# Row 132
$ws.Range("H132").Value = 3664.2307
$ws.Range("I132").Value = 3738.457
$ws.Range("J132").Value = 3014.75
$ws.Range("K132").Value = 11215.371
$ws.Range("L132").Value = 9044.25
$ws.Range("M132").Value = -8685.370999999999
$ws.Range("N132").Value = -14104.25
# Row 138
$ws.Range("H138").Value = 2424.6042
$ws.Range("I138").Value = 3497
$ws.Range("J138").Value = 2067.139
$ws.Range("K138").Value = 10491
$ws.Range("L138").Value = 6201.417
$ws.Range("M138").Value = -5351
$ws.Range("N138").Value = -16481.417

$ws = $wb.Worksheets.Item("ARM")
# Row 10
$ws.Range("H10").Value = 39000
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 39000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 39000
$ws.Range("M10").ClearContents() | Out-Null
$ws.Range("N10").Value = -39340
# Row 32
$ws.Range("H32").Value = 641318
$ws.Range("I32").Value = 832986.5600000001
$ws.Range("K32").Value = 832986.5600000001
$ws.Range("M32").Value = -832699.5600000001
# Row 45
$ws.Range("H45").Value = 2498.4348
$ws.Range("I45").Value = 1457.6923
$ws.Range("K45").Value = 1457.6923
$ws.Range("M45").Value = -1080.6923
# Row 61
$ws.Range("H61").Value = 2325.9443
$ws.Range("I61").Value = 1570.4348
$ws.Range("K61").Value = 1570.4348
$ws.Range("M61").Value = -1358.4348
# Row 74
$ws.Range("H74").Value = 1261.2222
$ws.Range("I74").Value = 1138.1818
$ws.Range("K74").Value = 1138.1818
$ws.Range("M74").Value = -264.1818000000001
# Row 77
$ws.Range("H77").Value = 1261.2222
$ws.Range("I77").Value = 1138.1818
$ws.Range("K77").Value = 5690.909000000001
$ws.Range("M77").Value = -1322.909000000001
# Row 132
$ws.Range("H132").Value = 3945.2646
$ws.Range("I132").Value = 3791.6316
$ws.Range("J132").Value = 4139.8667
$ws.Range("K132").Value = 11374.8948
$ws.Range("L132").Value = 12419.6001
$ws.Range("M132").Value = -8844.8948
$ws.Range("N132").Value = -17479.6001
# Row 136
$ws.Range("H136").Value = 2325.9443
$ws.Range("I136").Value = 1570.4348
$ws.Range("K136").Value = 4711.3044
$ws.Range("M136").Value = -2161.3044

$ws = $wb.Worksheets.Item("BSM")
# Row 12
$ws.Range("H12").Value = 200
$ws.Range("I12").Value = 200
$ws.Range("K12").Value = 200
$ws.Range("M12").Value = -32
# Row 107
$ws.Range("H107").Value = 2669.5557
$ws.Range("I107").Value = 2000
$ws.Range("J107").Value = 3004.3333
$ws.Range("K107").Value = 2000
$ws.Range("L107").Value = 3004.3333
$ws.Range("M107").Value = -80
$ws.Range("N107").Value = -6844.3333
# Row 134
$ws.Range("H134").Value = 2443.575
$ws.Range("I134").Value = 2122.5151
$ws.Range("K134").Value = 6367.5453
$ws.Range("M134").Value = -3832.5453

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 306.375
$ws.Range("J22").Value = 555.5
$ws.Range("L22").Value = 555.5
$ws.Range("N22").Value = -1255.5
# Row 31
$ws.Range("H31").Value = 3985.8408
$ws.Range("I31").Value = 926.75
$ws.Range("J31").Value = 17751.75
$ws.Range("K31").Value = 926.75
$ws.Range("L31").Value = 17751.75
$ws.Range("M31").Value = -631.75
$ws.Range("N31").Value = -18341.75
# Row 34
$ws.Range("H34").Value = 3985.8408
$ws.Range("I34").Value = 926.75
$ws.Range("J34").Value = 17751.75
$ws.Range("K34").Value = 926.75
$ws.Range("L34").Value = 17751.75
$ws.Range("M34").Value = -724.75
$ws.Range("N34").Value = -18155.75
# Row 38
$ws.Range("H38").Value = 8037
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents() | Out-Null
# Row 46
$ws.Range("H46").Value = 8037
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents() | Out-Null
# Row 58
$ws.Range("H58").Value = 1288.1428
$ws.Range("I58").Value = 715.25
$ws.Range("J58").Value = 2052
$ws.Range("K58").Value = 715.25
$ws.Range("L58").Value = 2052
$ws.Range("M58").Value = -512.25
$ws.Range("N58").Value = -2458
# Row 132
$ws.Range("H132").Value = 23811998
$ws.Range("I132").Value = 940.6667
$ws.Range("J132").Value = 41670292
$ws.Range("K132").Value = 2822.0001
$ws.Range("L132").Value = 125010876
$ws.Range("M132").Value = -292.0001000000002
$ws.Range("N132").Value = -125015936
# Row 134
$ws.Range("H134").Value = 1666.6471
$ws.Range("I134").Value = 1588.8667
$ws.Range("K134").Value = 4766.6001
$ws.Range("M134").Value = -2231.6001
# Row 136
$ws.Range("H136").Value = 1288.1428
$ws.Range("I136").Value = 715.25
$ws.Range("J136").Value = 2052
$ws.Range("K136").Value = 2145.75
$ws.Range("L136").Value = 6156
$ws.Range("M136").Value = 404.25
$ws.Range("N136").Value = -11256

$ws = $wb.Worksheets.Item("CUL")
# Row 93
$ws.Range("H93").Value = 5304.2856
$ws.Range("J93").Value = 5673.846
$ws.Range("L93").Value = 17021.538
$ws.Range("N93").Value = -20765.538
# Row 115
$ws.Range("H115").Value = 4494.5835
$ws.Range("I115").Value = 3294.6
$ws.Range("K115").Value = 9883.799999999999
$ws.Range("M115").Value = -8708.799999999999
# Row 137
$ws.Range("H137").Value = 16380.286
$ws.Range("I137").Value = 25682.25
$ws.Range("J137").Value = 3977.6667
$ws.Range("K137").Value = 77046.75
$ws.Range("L137").Value = 11933.0001
$ws.Range("M137").Value = -71946.75
$ws.Range("N137").Value = -22133.0001
# Row 140
$ws.Range("H140").Value = 1774.7693
$ws.Range("I140").Value = 858.55554
$ws.Range("J140").Value = 3836.25
$ws.Range("K140").Value = 2575.66662
$ws.Range("L140").Value = 11508.75
$ws.Range("M140").Value = 2604.33338
$ws.Range("N140").Value = -21868.75

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1774.3334
$ws.Range("I102").Value = 1529.2
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 1529.2
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = 92.79999999999995
$ws.Range("N102").Value = -6244
# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents() | Out-Null
# Row 126
$ws.Range("H126").Value = 3115.8823
$ws.Range("I126").Value = 2991.4
$ws.Range("J126").Value = 3293.7144
$ws.Range("K126").Value = 8974.200000000001
$ws.Range("L126").Value = 9881.143199999999
$ws.Range("M126").Value = -6504.200000000001
$ws.Range("N126").Value = -14821.1432
# Row 132
$ws.Range("H132").Value = 3537.4285
$ws.Range("I132").Value = 3402.4
$ws.Range("J132").Value = 3612.4443
$ws.Range("K132").Value = 10207.2
$ws.Range("L132").Value = 10837.3329
$ws.Range("M132").Value = -7677.200000000001
$ws.Range("N132").Value = -15897.3329

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 166669780
$ws.Range("I7").Value = 333334600
$ws.Range("J7").Value = 4936.6665
$ws.Range("K7").Value = 333334600
$ws.Range("L7").Value = 4936.6665
$ws.Range("M7").Value = -333334488
$ws.Range("N7").Value = -5160.6665
# Row 40
$ws.Range("H40").Value = 62505370
$ws.Range("I40").Value = 500001950
$ws.Range("J40").Value = 5857.143
$ws.Range("K40").Value = 500001950
$ws.Range("L40").Value = 5857.143
$ws.Range("M40").Value = -500001814
$ws.Range("N40").Value = -6129.143
# Row 126
$ws.Range("H126").Value = 166669780
$ws.Range("I126").Value = 333334600
$ws.Range("J126").Value = 4936.6665
$ws.Range("K126").Value = 1000003800
$ws.Range("L126").Value = 14809.9995
$ws.Range("M126").Value = -1000001330
$ws.Range("N126").Value = -19749.9995
# Row 132
$ws.Range("H132").Value = 3294
$ws.Range("I132").Value = 2124.476
$ws.Range("J132").Value = 5048.2856
$ws.Range("K132").Value = 6373.428
$ws.Range("L132").Value = 15144.8568
$ws.Range("M132").Value = -3843.428
$ws.Range("N132").Value = -20204.8568

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1260.6666
$ws.Range("I113").Value = 1494
$ws.Range("J113").Value = 654
$ws.Range("K113").Value = 4482
$ws.Range("L113").Value = 1962
$ws.Range("M113").Value = -2312
$ws.Range("N113").Value = -6302
# Row 132
$ws.Range("H132").Value = 5749915
$ws.Range("I132").Value = 3076.647
$ws.Range("K132").Value = 9229.940999999999
$ws.Range("M132").Value = -6699.940999999999
# Row 136
$ws.Range("H136").Value = 1894.0193
$ws.Range("I136").Value = 1569.9714
$ws.Range("J136").Value = 2561.1765
$ws.Range("K136").Value = 4709.914199999999
$ws.Range("L136").Value = 7683.529500000001
$ws.Range("M136").Value = -2159.914199999999
$ws.Range("N136").Value = -12783.5295
